$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), copying the format of the existing
# header cell G1 ("sum") so it reuses the same bold/bordered/centered
# style rather than minting a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column: 0/1 flag per row (row 6 = 1).
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 1
    7 = 0
    8 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
